# DRI_IOM_V7.xlsx edit:
#   "shifted a column to the right (replacing a blank column) in DRI_IOM_V7.xlsx"
#
# On the AMDR_hi sheet, column H was a blank spacer column and column I held
# the protein_g values. Move I's contents (header + all data rows) into H,
# leaving I empty again, and update the view/selection state to match.

$wb = $excel.ActiveWorkbook

# --- AMDR_hi sheet: move protein_g data from column I into column H ---
$wsAmdrHi = $wb.Worksheets.Item("AMDR_hi")

# Clear column I's formatting first. This makes the <cols> block collapse
# back to two ranges (3-8 and 10-16384) instead of keeping a dedicated
# entry for column 9, matching a column that is blank/default again.
$wsAmdrHi.Columns.Item(9).ClearFormats() | Out-Null

for ($r = 1; $r -le 23; $r++) {
  $srcCell = $wsAmdrHi.Cells.Item($r, 9)   # column I
  $dstCell = $wsAmdrHi.Cells.Item($r, 8)   # column H
  $v = $srcCell.Value()
  if ($v -ne $null) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4104) | Out-Null  # xlPasteAll: copy value + formatting
    $srcCell.Clear() | Out-Null              # fully clear the now-empty source cell
  }
}

# The selection on AMDR_hi moved from L14 to M3
$wsAmdrHi.Range("M3").Select() | Out-Null

# --- The active tab moved from AMDR_hi to RDAorAI_minerals ---
$wsMinerals = $wb.Worksheets.Item("RDAorAI_minerals")
$wsMinerals.Activate() | Out-Null
$wsMinerals.Range("R1").Select() | Out-Null
